$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Insert a new row at position 91 (everything from old row 91 onward
#    shifts down by one row, e.g. the old row 91 -> row 92, old row 94 (the
#    "Detection: Tube lens" header, with the hyperlinked C94 cell) -> row 95).
# ---------------------------------------------------------------------------
$ws.Rows.Item(91).Insert()

# ---------------------------------------------------------------------------
# 2. Populate the new row 91 with the new 20X Mitutoyo objective entry.
# ---------------------------------------------------------------------------
$ws.Range("A91").Value2 = "378-847"
$ws.Range("B91").Value2 = "Mitutoyo"
$ws.Range("C91").Value2 = "Mitutoyo G Plan APO 20X/t3,5"
$ws.Range("D91").Value2 = 0
$ws.Range("E91").Value2 = 3670
$ws.Range("F91").Formula = "=E91*D91"
$ws.Range("G91").Value2 = "Detection path"
$ws.Range("I91").Value2 = "Corrected for 3.5 mm of glass (n1.52). Thread is unusual, but fits SM1."

# Row height, matching the other "Detection: objectives" rows (e.g. row 87/90)
$ws.Rows.Item(91).RowHeight = 15.6

# Fonts: columns A, D, E, G, I use the "Arial 12 black" font (same as the
# other objective rows, e.g. row 87/90)
foreach ($col in @("A91","D91","E91","G91","I91")) {
    $rng = $ws.Range($col)
    $rng.Font.Name = "Arial"
    $rng.Font.Size = 12
    $rng.Font.Color = 0
}

# Columns B and H use the plain "Arial 12" font (no explicit color)
foreach ($col in @("B91","H91")) {
    $rng = $ws.Range($col)
    $rng.Font.Name = "Arial"
    $rng.Font.Size = 12
}

# E91 / F91 are price/total cells -> 2 decimal number format, like the other
# price cells in this section
$ws.Range("E91").NumberFormat = "0.00"
$ws.Range("F91").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# 3. The MT-1 tube lens note (row 89) becomes "Optional" and its quantity
#    drops to 0.
# ---------------------------------------------------------------------------
$ws.Range("D89").Value2 = 0
$ws.Range("I89").Value2 = "Works best with the Mitutoyo MT-1 tube lens. Optional."

# ---------------------------------------------------------------------------
# 4. The black-ring back cover (old row 91, now row 92) needs 5 instead of 4.
# ---------------------------------------------------------------------------
$ws.Range("D92").Value2 = 5

# ---------------------------------------------------------------------------
# 5. Fix up hyperlinks. The row insert does not automatically repoint the
#    hyperlink that used to sit on C94 (the MT-1 tube lens row), which is now
#    C95, so every hyperlink on the sheet is recreated pointing at the
#    correct (possibly shifted) cells, plus the brand-new one for C91.
# ---------------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C87"), "https://www.edmundoptics.com/p/2x-mitutoyo-bd-plan-apo-objective/45425/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C88"), "https://www.edmundoptics.com/p/5x-mitutoyo-bd-plan-apo-objective/45426/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C89"), "https://www.edmundoptics.com/p/75x-mitutoyo-bd-plan-apo-objective/45427/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C90"), "https://www.edmundoptics.com/p/10x-mitutoyo-bd-plan-apo-objective/45428/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C95"), "https://www.edmundoptics.com/p/mt-1-accessory-tube-lens/11488/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C78"), "https://astronomy-imaging-camera.com/product/efw-mini") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C91"), "https://www.mitutoyo.com/webfoo/wp-content/uploads/Objectives_Catalog_E13014.pdf") | Out-Null
